# Applies the "cryptos" price/volume update described in the commit diff.
# Column layout: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.806.40'
$ws.Range("E2").Value = '  +5.19%  '

$ws.Range("D3").Value = '2.275.72'
$ws.Range("E3").Value = '  +3.25%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '''233.41'
$ws.Range("E5").Value = '  +1.65%  '

$ws.Range("D6").Value = '''0.639'
$ws.Range("E6").Value = '  +4.32%  '

$ws.Range("D7").Value = '''66.04'
$ws.Range("E7").Value = '  +9.54%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("E9").Value = '  +6.72%  '

$ws.Range("D10").Value = '''0.104'
$ws.Range("E10").Value = '  +16.80%  '

$ws.Range("D11").Value = '''57.56'
$ws.Range("E11").Value = '  +0.82%  '

$ws.Range("D12").Value = '''26.20'
$ws.Range("E12").Value = '  +17.94%  '

$ws.Range("E13").Value = '  +0.72%  '

$ws.Range("D14").Value = '2.621.09'
$ws.Range("E14").Value = '  +3.59%  '

$ws.Range("D15").Value = '''15.96'
$ws.Range("E15").Value = '  +3.85%  '

$ws.Range("D16").Value = '''5.97'
$ws.Range("E16").Value = '  +5.65%  '

$ws.Range("D17").Value = '''0.833'
$ws.Range("E17").Value = '  +5.09%  '

$ws.Range("D18").Value = '2.277.07'
$ws.Range("E18").Value = '  +3.56%  '

$ws.Range("D19").Value = '43.714.09'
$ws.Range("E19").Value = '  +5.11%  '

$ws.Range("D20").Value = '0.0₃0995'
$ws.Range("E20").Value = '  +10.55%  '

$ws.Range("D21").Value = '''74.60'
$ws.Range("E21").Value = '  +3.57%  '

$ws.Range("D22").Value = '''6.16'
$ws.Range("E22").Value = '  +1.88%  '

$ws.Range("D23").Value = '''262.80'
$ws.Range("E23").Value = '  +8.39%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("E25").Value = '  +5.97%  '

$ws.Range("D26").Value = '''2.33'
$ws.Range("E26").Value = '  +1.91%  '

$ws.Range("D27").Value = '''10.10'
$ws.Range("E27").Value = '  +5.02%  '

$ws.Range("D28").Value = '''173.32'
$ws.Range("E28").Value = '  +2.45%  '

$ws.Range("D29").Value = '''21.15'
$ws.Range("E29").Value = '  +7.17%  '

$ws.Range("D30").Value = '''0.137'
$ws.Range("E30").Value = '  -1.88%  '

$ws.Range("E31").Value = '  -0.68%  '

$ws.Range("D32").Value = '''2.80'
$ws.Range("E32").Value = '  +8.61%  '

$ws.Range("E33").Value = '  +3.07%  '

$ws.Range("D34").Value = '''0.0689'
$ws.Range("E34").Value = '  +6.64%  '

$ws.Range("D35").Value = '''5.09'
$ws.Range("E35").Value = '  +1.55%  '

$ws.Range("E36").Value = '  +3.54%  '

$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '''6.78'
$ws.Range("E37").Value = '  +7.11%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''3.82'
$ws.Range("E38").Value = '  +8.12%  '

$ws.Range("E39").Value = '  +0.57%  '

$ws.Range("E40").Value = '  +5.14%  '

$ws.Range("E41").Value = '  +0.15%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''8.44'
$ws.Range("E42").Value = '  -1.01%  '

$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").Value = '''4.59'
$ws.Range("E43").Value = '  +4.83%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''17.59'
$ws.Range("E44").Value = '  +7.62%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '''0.0981'
$ws.Range("E45").Value = '  +3.36%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '''10.55'
$ws.Range("E46").Value = '  +23.01%  '

$ws.Range("D47").Value = '''98.66'
$ws.Range("E47").Value = '  +1.45%  '

$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''2.39'
$ws.Range("E49").Value = '  +7.93%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.483.80'
$ws.Range("E50").Value = '  +1.26%  '

$ws.Range("D51").Value = '''0.000209'
$ws.Range("E51").Value = '  -12.45%  '
